$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.053.59'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.679.68'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.84'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.253'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.09%  '
$ws.Range('E9').Value = '  +5.81%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0623'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0890'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '1.917.58'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '1.675.86'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.533'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.78%  '
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').Value = '27.061.17'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '236.47'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.46'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.59%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.26'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.66%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.21'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('E27').Value = '  +3.90%  '
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.36'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').Value = '1.544.14'
$ws.Range('E33').Value = '  +6.09%  '
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range('E35').Value = '  +5.49%  '
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.589'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.58%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.914'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('E39').Value = '  +3.06%  '
$ws.Range('E40').Value = '  +6.75%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '67.85'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.53'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.25'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('D45').Value = '1.822.57'
$ws.Range('E45').Value = '  +0.92%  '
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.52'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').Value = '  +2.19%  '
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.04'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +6.54%  '
